$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5: new candidate IDs / credentials ---
$ws.Range("A2").Value = "ugslZ418"
$ws.Range("B2").Value = 231011196
$ws.Range("C2").Value = "ltottxv52"
$ws.Range("D2").Value = "Q&x!67Sq"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "UsyZihDX"
$ws.Range("G2").Value = "xvsd"
$ws.Range("H2").Value = "Candidate"

$ws.Range("A3").Value = "EymBh619"
$ws.Range("B3").Value = 231011195
$ws.Range("C3").Value = "ezvsenq19"
$ws.Range("D3").Value = "n%5hPS4!"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "jSOQfLlV"
$ws.Range("G3").Value = "WrSW"
$ws.Range("H3").Value = "Candidate"

$ws.Range("A4").Value = "MPszX650"
$ws.Range("B4").Value = 231011194
$ws.Range("C4").Value = "ftqbujt18"
$ws.Range("D4").Value = "D8`$dSq3!"
$ws.Range("E4").Value = "MR"
$ws.Range("F4").Value = "NvkYrwyA"
$ws.Range("G4").Value = "FaXr"
$ws.Range("H4").Value = "Candidate"

$ws.Range("A5").Value = "VFUmf147"
$ws.Range("B5").Value = 231011193
$ws.Range("C5").Value = "ubaytlu22"
$ws.Range("D5").Value = "b!DM4r3&"
$ws.Range("E5").Value = "MR"
$ws.Range("F5").Value = "MQwRDcpH"
$ws.Range("G5").Value = "REwv"
$ws.Range("H5").Value = "Candidate"

# --- Add new row 6, carrying the same look & feel (borders/format) as row 5 ---
$ws.Range("A5:H5").Copy($ws.Range("A6:H6"))

$ws.Range("A6").Value = "mZIYS782"
$ws.Range("B6").Value = 231011192
$ws.Range("C6").Value = "hyamvsw72"
$ws.Range("D6").Value = "EZk62w`$#"
$ws.Range("E6").Value = "MR"
$ws.Range("F6").Value = "xLLKASPY"
$ws.Range("G6").Value = "ivpn"
$ws.Range("H6").Value = "Candidate"

# --- Keep dimension / selection in sync with the now-6-row table ---
$ws.Range("A1:H6").Select() | Out-Null
